# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last updated" timestamp header ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 17:35"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1710558
$ws.Range("C4").Value = 4332
$ws.Range("D4").Value = 466980
$ws.Range("E4").Value = 1143650
$ws.Range("G4").Value = 123
$ws.Range("H4").Value = 99928

# --- Row 13: India ---
$ws.Range("B13").Value = 149822
$ws.Range("C13").Value = 4872
$ws.Range("D13").Value = 63465
$ws.Range("E13").Value = 82050
$ws.Range("G13").Value = 135
$ws.Range("H13").Value = 4307

# --- Rows 18/19: Chile overtakes Arabia Saudita, rows swap places ---
$ws.Range("A18").Value = "Chile"
$ws.Range("B18").Value = 77961
$ws.Range("C18").Value = 3964
$ws.Range("D18").Value = 30915
$ws.Range("E18").Value = 46240
$ws.Range("G18").Value = 45
$ws.Range("H18").Value = 806

$ws.Range("A19").Value = "Arabia Saudita"
$ws.Range("B19").Value = 76726
$ws.Range("C19").Value = 1931
$ws.Range("D19").Value = 48450
$ws.Range("E19").Value = 27865
$ws.Range("G19").Value = 12
$ws.Range("H19").Value = 411

# --- Row 29: Singapur ---
$ws.Range("D29").Value = 16444
$ws.Range("E29").Value = 15876

# --- Row 76: Uzbekistan ---
$ws.Range("D76").Value = 2636
$ws.Range("E76").Value = 632

# --- Row 81: Grecia ---
$ws.Range("B81").Value = 2892
$ws.Range("C81").Value = 10
$ws.Range("E81").Value = 1345
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = 173

# --- Row 94: Somalia ---
$ws.Range("B94").Value = 1711
$ws.Range("C94").Value = 22
$ws.Range("D94").Value = 253
$ws.Range("E94").Value = 1391
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 67

# --- Row 116: Republica de Chipre ---
$ws.Range("B116").Value = 939
$ws.Range("C116").Value = 2
$ws.Range("E116").Value = 328

# --- Rows 126/127: Jordania overtakes Crucero, rows swap places ---
$ws.Range("A126").Value = "Jordania"
$ws.Range("B126").Value = 718
$ws.Range("C126").Value = 7
$ws.Range("D126").Value = 479
$ws.Range("E126").Value = 230
$ws.Range("H126").Value = 9

$ws.Range("A127").Value = "Crucero"
$ws.Range("B127").Value = 712
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 651
$ws.Range("E127").Value = 48
$ws.Range("H127").Value = 13

# --- Row 137: Reunion ---
$ws.Range("B137").Value = 459
$ws.Range("C137").Value = 3
$ws.Range("D137").Value = 47
$ws.Range("E137").Value = 47
